$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("D1").Value = "porcentaje_utilidades"
$ws.Range("E1").Value = "porcentaje_contingencia"

# Add new data values
$ws.Range("E2").Value = 13
$ws.Range("D3").Value = 11.11

# Adjust column widths for new columns D and E
$ws.Columns.Item(4).ColumnWidth = 19.90625
$ws.Columns.Item(5).ColumnWidth = 22.36328125

# Update selection to E4 as in the target sheetView
$ws.Range("E4").Select()
